$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates scraped from the commit diff: (cellRef, newValue).
$updates = @(
    @('D2', '64.158.12'),
    @('E2', '  +1.01%  '),
    @('D3', '3.124.70'),
    @('E3', '  +1.28%  '),
    @('E4', '  +0.01%  '),
    @('D5', '602.19'),
    @('E5', '  -0.67%  '),
    @('D6', '142.44'),
    @('E6', '  -0.99%  '),
    @('E7', '  -0.10%  '),
    @('D8', '3.122.58'),
    @('E8', '  +1.29%  '),
    @('D9', '0.522'),
    @('E9', '  +0.93%  '),
    @('E10', '  +1.04%  '),
    @('E11', '  +3.55%  '),
    @('D12', '0.467'),
    @('E12', '  +0.48%  '),
    @('E13', '  +3.88%  '),
    @('D14', '35.10'),
    @('E14', '  +0.59%  '),
    @('D15', '3.640.99'),
    @('E15', '  +1.15%  '),
    @('E16', '  +3.15%  '),
    @('D17', '64.019.48'),
    @('E17', '  +0.62%  '),
    @('D18', '3.117.34'),
    @('E18', '  +0.94%  '),
    @('D19', '6.86'),
    @('E19', '  +1.64%  '),
    @('D20', '477.69'),
    @('E20', '  +1.13%  '),
    @('E21', '  +0.19%  '),
    @('D22', '0.710'),
    @('E22', '  +1.83%  '),
    @('D23', '7.66'),
    @('E23', '  +0.50%  '),
    @('D24', '85.17'),
    @('E24', '  +2.68%  '),
    @('D25', '13.34'),
    @('E25', '  -0.59%  '),
    @('D26', '1.00'),
    @('E26', '  -0.07%  '),
    @('E27', '  -0.29%  '),
    @('D28', '8.33'),
    @('E28', '  +0.64%  '),
    @('D29', '7.16'),
    @('E29', '  +8.15%  '),
    @('B30', 'ImmutableX'),
    @('C30', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'),
    @('D30', '2.04'),
    @('E30', '  -3.92%  '),
    @('B31', 'Hedera'),
    @('C31', 'https://coinranking.com/coin/jad286TjB+hedera-hbar'),
    @('D31', '0.113'),
    @('E31', '  +1.46%  '),
    @('E32', '  -0.03%  '),
    @('D33', '26.78'),
    @('E33', '  +3.29%  '),
    @('E34', '  -2.38%  '),
    @('E35', '  +0.02%  '),
    @('B36', 'Filecoin'),
    @('C36', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'),
    @('D36', '5.95'),
    @('E36', '  +1.28%  '),
    @('B37', 'PEPE'),
    @('C37', 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'),
    @('D37', '0.0₃0765'),
    @('E37', '  +6.06%  '),
    @('E38', '  +0.02%  '),
    @('E39', '  +4.47%  '),
    @('D40', '443.25'),
    @('E40', '  -1.66%  '),
    @('E41', '  +0.68%  '),
    @('E42', '  +0.90%  '),
    @('D43', '8.19'),
    @('E43', '  -0.99%  '),
    @('D44', '2.851.57'),
    @('E44', '  +1.56%  '),
    @('E45', '  -1.22%  '),
    @('E46', '  +0.06%  '),
    @('D47', '2.42'),
    @('E47', '  +2.04%  '),
    @('E48', '  +0.02%  '),
    @('D49', '25.93'),
    @('E49', '  +0.50%  '),
    @('E50', '  +0.49%  '),
    @('D51', '120.16'),
    @('E51', '  +2.44%  '),
)

foreach ($update in $updates) {
    $cellRef = $update[0]
    $text = $update[1]
    # Plain numeric-looking text (e.g. crypto "Price" values like "602.19")
    # would otherwise be auto-coerced to a real number by Range.Value (and
    # lose things like trailing zeros, e.g. "1.00" -> 1). A leading apostrophe
    # keeps it text, exactly like a user typing '602.19 into the cell.
    if ($text -match '^[0-9]+(\.[0-9]+)?$') {
        $text = "'" + $text
    }
    $ws.Range($cellRef).Value = $text
}
